$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.834.03"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.094.22"
$ws.Range("E3").Value = "  -2.87%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.76"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.32"
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.540"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.093.08"
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("E12").Value = "  -4.07%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -5.69%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.95"
$ws.Range("E14").Value = "  -5.79%  "
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.603.25"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.22"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.741.55"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.095.29"
$ws.Range("E19").Value = "  -3.10%  "
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.46"
$ws.Range("E21").Value = "  -3.94%  "
$ws.Range("E22").Value = "  -7.29%  "
$ws.Range("E23").Value = "  -3.39%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.23"
$ws.Range("E25").Value = "  -2.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.91"
$ws.Range("E26").Value = "  -5.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.63"
$ws.Range("E27").Value = "  +6.26%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.52"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.18"
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("E33").Value = "  -5.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.23"
$ws.Range("E34").Value = "  -4.02%  "
$ws.Range("E35").Value = "  -6.69%  "
$ws.Range("E36").Value = "  -2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.02"
$ws.Range("E37").Value = "  -4.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.25"
$ws.Range("E38").Value = "  -4.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("E39").Value = "  -8.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.93"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.24"
$ws.Range("E41").Value = "  -2.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "439.75"
$ws.Range("E42").Value = "  -6.98%  "
$ws.Range("E43").Value = "  -4.19%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0363"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.112"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.13"
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.825.85"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.18"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  -3.99%  "
